# Anonymize SampleBook1: replace First Name / Last Name columns on Sheet1
# with the (anonymized) STDID number, and replace the STDID values
# themselves (on both sheets) with new anonymized numeric ids.
# Header rows / other columns are left untouched.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Mapping from the original STDID to the new anonymized STDID.
$idMap = @{
    12345678 = 674920753
    23456789 = 406378728
    34567890 = 852158896
    10101010 = 123456789
    23674567 = 349854278
}

# --- Sheet1 ---------------------------------------------------------------
# Columns: A=STDID  B=First Name  C=Last Name  D=Gender  E=Secion  F=Grade
#          G=Prereq?  H=Comment
for ($r = 2; $r -le 6; $r++) {
    $oldId = [int]$ws1.Cells.Item($r, 1).Value2
    $newId = $idMap[$oldId]

    $ws1.Cells.Item($r, 1).Value = $newId   # STDID
    $ws1.Cells.Item($r, 2).Value = $newId   # First Name -> anonymized id
    $ws1.Cells.Item($r, 3).Value = $newId   # Last Name  -> anonymized id
}

# --- Sheet2 ---------------------------------------------------------------
# Columns: A=STDID  B=Didtest  C=TestScore
for ($r = 2; $r -le 5; $r++) {
    $oldId = [int]$ws2.Cells.Item($r, 1).Value2
    $newId = $idMap[$oldId]

    $ws2.Cells.Item($r, 1).Value = $newId   # STDID
}
